$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---- Row 8: ECD-I | Knime | Program | https://www.knime.com/ ----
$ws.Cells.Item(8, 1).Value = "ECD-I"
$ws.Cells.Item(8, 4).Value = "https://www.knime.com/"
$ws.Cells.Item(8, 2).Value = "Knime"
$ws.Cells.Item(8, 3).Value = "Program"

$d8 = $ws.Cells.Item(8, 4)
$ws.Hyperlinks.Add($d8, "https://www.knime.com/") | Out-Null

# Re-apply the same direct formatting used by the existing rows so the
# new cells land on the same (already-present) cell styles rather than
# picking up ad-hoc ones.
$ws.Range("C7").Copy() | Out-Null
$ws.Cells.Item(8, 3).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D7").Copy() | Out-Null
$ws.Cells.Item(8, 4).PasteSpecial($xlPasteFormats) | Out-Null

# ---- Row 9: ECD-I | Weka | Program | http://www.cs.waikato.ac.nz/ml/weka/ ----
$ws.Cells.Item(9, 1).Value = "ECD-I"
$ws.Cells.Item(9, 2).Value = "Weka"
$ws.Cells.Item(9, 3).Value = "Program"
$ws.Cells.Item(9, 4).Value = "http://www.cs.waikato.ac.nz/ml/weka/"

$d9 = $ws.Cells.Item(9, 4)
$ws.Hyperlinks.Add($d9, "http://www.cs.waikato.ac.nz/ml/weka/") | Out-Null

$ws.Range("D7").Copy() | Out-Null
$ws.Cells.Item(9, 4).PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# Matches the post-edit cursor position recorded in the workbook.
$ws.Range("D10").Select() | Out-Null
